$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "62.047.32"
Set-TextValue "E2" "  +2.37%  "

Set-TextValue "D3" "2.407.40"
Set-TextValue "E3" "  -0.30%  "

Set-TextValue "E4" "  +0.81%  "

Set-TextValue "D5" "570.93"
Set-TextValue "E5" "  +1.13%  "

Set-TextValue "D6" "144.12"
Set-TextValue "E6" "  +4.61%  "

Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.50%  "

Set-TextValue "D8" "0.537"
Set-TextValue "E8" "  +0.58%  "

Set-TextValue "D9" "2.429.70"
Set-TextValue "E9" "  +1.28%  "

Set-TextValue "E10" "  +4.28%  "

Set-TextValue "E11" "  +0.79%  "

Set-TextValue "D12" "5.24"
Set-TextValue "E12" "  +4.06%  "

Set-TextValue "E13" "  +3.56%  "

Set-TextValue "D14" "26.63"
Set-TextValue "E14" "  +3.64%  "

Set-TextValue "D15" "0.0000177"
Set-TextValue "E15" "  +6.18%  "

Set-TextValue "E16" "  +1.64%  "

Set-TextValue "D17" "61.857.49"
Set-TextValue "E17" "  +2.05%  "

Set-TextValue "D18" "2.427.00"
Set-TextValue "E18" "  +1.30%  "

Set-TextValue "D19" "7.92"
Set-TextValue "E19" "  -3.45%  "

Set-TextValue "D20" "10.82"
Set-TextValue "E20" "  +2.41%  "

Set-TextValue "D21" "325.40"
Set-TextValue "E21" "  +0.85%  "

Set-TextValue "E22" "  +2.28%  "

Set-TextValue "E23" "  +12.79%  "

Set-TextValue "E24" "  -0.14%  "

Set-TextValue "D25" "65.21"
Set-TextValue "E25" "  +1.58%  "

Set-TextValue "D26" "613.02"
Set-TextValue "E26" "  +10.91%  "

Set-TextValue "D27" "8.41"
Set-TextValue "E27" "  +4.53%  "

Set-TextValue "D28" "0.0₃0981"
Set-TextValue "E28" "  +7.86%  "

Set-TextValue "D30" "8.05"
Set-TextValue "E30" "  +2.50%  "

Set-TextValue "E31" "  +8.82%  "

Set-TextValue "E32" "  +1.72%  "

Set-TextValue "E33" "  +3.10%  "

Set-TextValue "E34" "  +4.27%  "

Set-TextValue "D35" "0.996"
Set-TextValue "E35" "  -0.71%  "

Set-TextValue "E36" "  +5.77%  "

Set-TextValue "E37" "  +0.23%  "

Set-TextValue "D38" "0.372"
Set-TextValue "E38" "  +1.25%  "

Set-TextValue "E39" "  +5.86%  "

Set-TextValue "E40" "  +1.59%  "

Set-TextValue "D41" "2.66"
Set-TextValue "E41" "  +16.78%  "

Set-TextValue "D42" "1.73"
Set-TextValue "E42" "  +5.32%  "

Set-TextValue "E43" "  -0.06%  "

Set-TextValue "D44" "42.13"
Set-TextValue "E44" "  +1.62%  "

Set-TextValue "E45" "  -2.40%  "

Set-TextValue "D46" "142.78"
Set-TextValue "E46" "  +0.27%  "

Set-TextValue "E47" "  +2.51%  "

Set-TextValue "D48" "20.25"
Set-TextValue "E48" "  +6.58%  "

Set-TextValue "E49" "  +2.65%  "

Set-TextValue "D50" "0.0511"
Set-TextValue "E50" "  +3.26%  "

Set-TextValue "D51" "0.0917"
Set-TextValue "E51" "  +2.67%  "

Write-Host "Done"